$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 3) to the sheet, duplicating the data row above it
# (row 2), per the commit's Katalon-AI-generated selector/option data set.
$ws.Range("A3").Value = 'products__item\ in-stock\ products__item_3-in-row"]:nth-child(1) [type="button'
$ws.Range("B3").Value = 'Hot! New!iPhone 15From:$'
$ws.Range("C3").Value = '\31 52170-case-636'
$ws.Range("D3").Value = '\31 52171-case-641'
$ws.Range("E3").Value = '1 TB'
$ws.Range("F3").Value = 'Black'
